$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "SamplesTab" query (B3) is rewritten: the extra "Tumor" /
# "Analyte Type" columns are removed from the SELECT list and the
# FROM clause is folded onto the Accession line.
$ws.Range("B3").Value = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
     s.phs_accession = 'phs001819' AND f1.experimental_strategy_and_data_subtypes = 'WXS'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

# Reflect the author's final cursor position / scroll state: the
# window had scrolled down/left so row 3 is the top row, with C3 as
# the active selected cell.
$win = $excel.ActiveWindow
$win.ScrollRow = 3
$win.ScrollColumn = 1
$ws.Range("C3").Select()
